# Added the batch integration tests
#
# The commit adds "Yes" answers for the newly-covered batch endpoints
# (rows 11-15, column C) and leaves the cursor/selection on C16 (instead
# of the previous scrolled-down B2 selection).  A couple of nearby cells
# (A22/A25, B22/B25) also pick up the sheet's normal (border-less) style
# instead of a stray duplicate style that Excel later garbage-collected
# when it rewrote styles.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- content: batch integration tests now covered ---
$ws.Range("C11").Value = "Yes"
$ws.Range("C12").Value = "Yes"
$ws.Range("C13").Value = "Yes"
$ws.Range("C14").Value = "Yes"
$ws.Range("C15").Value = "Yes"

# --- normalize a couple of stray cell styles to match their neighbours ---
$ws.Range("A22").Borders.LineStyle = -4142
$ws.Range("A25").Borders.LineStyle = -4142
$ws.Range("B22").Borders.LineStyle = -4142
$ws.Range("B25").Borders.LineStyle = -4142

# --- selection / scroll position ends up on C16 ---
$ws.Range("C16").Select()
